$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TO DO LIST")

# Update the "Projected Completion" dates in the lower rows of the TO DO LIST
# sheet (serial date numbers, keeping existing cell formatting/style).
$ws.Range("D17").Value = 43795
$ws.Range("D18").Value = 43797
$ws.Range("D19").Value = 43797
$ws.Range("D20").Value = 43797
$ws.Range("D21").Value = 43797
